$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new label/comment cell value in C40
$ws.Range("C40").Value = "option for cNORM recoding, express grade as weeks in school"

# Add the weeks-in-school values for the grade strata rows (41-46), column C
$ws.Range("C41").Value = 8
$ws.Range("C42").Value = 34
$ws.Range("C43").Value = 60
$ws.Range("C44").Value = 86
$ws.Range("C45").Value = 112
$ws.Range("C46").Value = 138

# Restore the view state: scroll to A22 and select C41
$ws.Range("C41").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
